$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("新题")

# The "排序" label in A3 was a stray leftover from an earlier row layout;
# it is no longer referenced anywhere else, so just clear it out.
$ws.Range("A3").ClearContents()

# Log a new day's entry: 2019-03-20 (serial 43544) -> "70 recursion", done.
# Copy/PasteSpecial(formats) from the previous date cell so the new cell
# reuses the existing date-number-format style instead of creating a new one.
$ws.Range("A11").Value = 43544
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Value = "70 recursion"
$ws.Range("E11").Value = "done"

# The longer date column now needs to widen to fit its content, same as
# column B already does (bestFit).
$ws.Range("A:A").AutoFit() | Out-Null

# Leave the selection where the author left it when they saved.
$ws.Range("D17:D18").Select() | Out-Null
